# Adds two new sheets (Sheet2, Sheet3) with product inventory test data,
# mirroring the data-driven test fixtures added for the Product APIs.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Sheet2: "raw" product rows backed by a real Excel Table (Table2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$headers2 = @("name","sku","description","initialStock","currentStock","reorderLevel","costPrice","Category","reservedStock")
for ($i = 0; $i -lt $headers2.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers2[$i]
}

$rows2 = @(
    @("Product1001","SKUTEST1001","Test Description for Product 1001",20,20,5,100,"TestCategory",0),
    @("Product1002","SKUTEST1002","Test Description for Product 1002",20,20,5,100,"TestCategory",0),
    @("Product1003","SKUTEST1003","Test Description for Product 1003",20,20,5,100,"TestCategory",0)
)
for ($r = 0; $r -lt $rows2.Length; $r++) {
    $row = $rows2[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# costPrice column (G) is formatted with 2 decimal places
$ws2.Range("G2:G4").NumberFormat = "0.00"

$tbl2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:I4"), $null, 1)
$tbl2.Name = "Table2"
$tbl2.TableStyle = "TableStyleMedium2"

$ws2.Range("A1:H4").Select()

# ---------------------------------------------------------------------------
# Sheet3: "changed" product rows (updated sku/description/category values)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

$headers3 = @("name","sku","description","reorderLevel","costPrice","Category")
for ($i = 0; $i -lt $headers3.Length; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $headers3[$i]
}

$rows3 = @(
    @("Product1001","SKUTEST#1001","Changed Test Description for Product 1001",15,120.5,"TestCategory1"),
    @("Product1002","SKUTEST#1002","Changed Test Description for Product 1002",20,104.5,"TestCategory2"),
    @("Product1003","SKUTEST#1003","Changed Test Description for Product 1003",10,90,"TestCategory3")
)
for ($r = 0; $r -lt $rows3.Length; $r++) {
    $row = $rows3[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws3.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# costPrice column (E) is formatted with 2 decimal places
$ws3.Range("E2:E4").NumberFormat = "0.00"

$ws3.Range("F5").Select()

Write-Output "done"
